$pairs = @(
    @("2025-06-19 Thursday", "2025-06-20 Friday"),
    @("73-17=", "50-26="),
    @("64-9=", "32-20="),
    @("83-15=", "74-14="),
    @("17+69=", "17+28="),
    @("4+79=", "49-35="),
    @("88-25=", "42+31="),
    @("74-44=", "60+17="),
    @("74+5=", "62+29="),
    @("42+15=", "29+65="),
    @("63+16=", "53+37="),
    @("52-27=", "81-35="),
    @("80-62=", "33-22="),
    @("20+49=", "48+47="),
    @("75-55=", "55-19="),
    @("11-2=", "65-15="),
    @("46+13=", "34-2="),
    @("73-18=", "0+80="),
    @("94-35=", "43-32="),
    @("59-56=", "79+9="),
    @("70-17=", "29+38="),
    @("85+6=", "7+88="),
    @("32+39=", "44+48="),
    @("87-64=", "69-26="),
    @("53+40=", "48+44="),
    @("18+54=", "5+51="),
    @("3+9=", "35+64="),
    @("16+12=", "52+27="),
    @("17-11=", "0+69="),
    @("55-37=", "34-29="),
    @("28+12=", "20-18="),
    @("47-31=", "42-11="),
    @("99-72=", "54-51="),
    @("75-2=", "91-70="),
    @("51-42=", "40-39="),
    @("88-2=", "89-25="),
    @("29-24=", "35-25="),
    @("44+31=", "85-65="),
    @("22+36=", "72-39="),
    @("95-62=", "95-27="),
    @("90-63=", "20+68="),
    @("30+24=", "17+71="),
    @("96-54=", "9-5="),
    @("69-38=", "82-78="),
    @("77-42=", "73-9="),
    @("16+13=", "48+47="),
    @("66+6=", "23+15="),
    @("84-47=", "73+8="),
    @("61+24=", "0+38="),
    @("18+41=", "44+8="),
    @("33+2=", "36-32="),
    @("1+93=", "10+33="),
    @("37+30=", "41-4="),
    @("93-48=", "47+12="),
    @("93-86=", "25-13="),
    @("27+70=", "56+41="),
    @("21+73=", "79+4="),
    @("79-21=", "25-15="),
    @("21-9=", "53-22="),
    @("8+68=", "60-4="),
    @("27+15=", "30-13="),
    @("88-20=", "13+58="),
    @("3+7=", "81-2="),
    @("45+31=", "74-65="),
    @("17-1=", "19+63="),
    @("5+4=", "18+47="),
    @("52+0=", "35+56="),
    @("22+26=", "8+46="),
    @("16+73=", "26-12="),
    @("59-7=", "9+60="),
    @("44+51=", "58+27="),
    @("98-7=", "13-10="),
    @("7-5=", "69+17="),
    @("25-5=", "69-3="),
    @("48-39=", "60+13="),
    @("88+4=", "40-8="),
    @("59-52=", "73+21="),
    @("63+25=", "24+30="),
    @("41+11=", "69+0="),
    @("43+28=", "21-13="),
    @("2+24=", "53+24="),
    @("83-21=", "69+2="),
    @("77-34=", "73+9="),
    @("92-36=", "84-81="),
    @("48+16=", "94-28="),
    @("98-96=", "72+4="),
    @("81-44=", "30+34="),
    @("49+5=", "14+25="),
    @("48-22=", "28+29="),
    @("56-55=", "92-3="),
    @("0+23=", "28-7="),
    @("72-55=", "83-2="),
    @("36+8=", "16-16="),
    @("44+5=", "61-46="),
    @("51-10=", "65+13="),
    @("77-45=", "1+61="),
    @("44-5=", "70-34="),
    @("62+26=", "70-19="),
    @("25-10=", "58+7="),
    @("56+24=", "34-2="),
    @("36-33=", "33-23="),
)
$d = $word.ActiveDocument

foreach ($pair in $pairs) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $oldText"
    }
}
